$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the student's name (typo "HURTADO HURTADO" -> correct full name)
$ws.Range("C3").Value = "BANEGAS HURTADO KAREN KATHERINE"

# Record today's (day 4 / column H) attendance for every student
$ws.Range("H3").Value = "p"
$ws.Range("H4").Value = "p"
$ws.Range("H5").Value = "j"
$ws.Range("H6").Value = "p"
$ws.Range("H7").Value = "p"
$ws.Range("H8").Value = "p"
$ws.Range("H9").Value = "p"
$ws.Range("H10").Value = "P"
$ws.Range("H11").Value = "P"
$ws.Range("H12").Value = "P"
$ws.Range("H13").Value = "P"
$ws.Range("H14").Value = "P"
$ws.Range("H15").Value = "p"
$ws.Range("H16").Value = "J"
$ws.Range("H17").Value = "p"
$ws.Range("H18").Value = "p"
$ws.Range("H19").Value = "p"
$ws.Range("H21").Value = "p"
$ws.Range("H22").Value = "p"

# Update the active view to reflect today's column being worked on
$ws.Range("H6").Select()
